# Apply cryptos list update (prices + 1h volume changes), Sat Sep 14 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force a literal TEXT value (no number/date auto-coercion, no style change)
    # by writing it as a string-literal formula, then collapsing the formula to a
    # static value via copy / paste-special-values (xlPasteValues = -4163).
    $escaped = $text.Replace('"', '""')
    $cell = $ws.Range($cellRef)
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2
Set-TextValue "D2" "60.009.58"
$ws.Range("E2").Value = "  +3.34%  "

# Row 3
Set-TextValue "D3" "2.423.09"
$ws.Range("E3").Value = "  +2.86%  "

# Row 4
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
Set-TextValue "D5" "553.16"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6
Set-TextValue "D6" "138.02"
$ws.Range("E6").Value = "  +3.99%  "

# Row 7
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
Set-TextValue "D8" "0.588"
$ws.Range("E8").Value = "  +3.29%  "

# Row 9
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
Set-TextValue "D10" "5.69"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D11" "0.356"
$ws.Range("E11").Value = "  +0.51%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.148"
$ws.Range("E12").Value = "  -2.14%  "

# Row 13
Set-TextValue "D13" "25.20"
$ws.Range("E13").Value = "  +4.63%  "

# Row 14
Set-TextValue "D14" "2.854.24"
$ws.Range("E14").Value = "  +2.61%  "

# Row 15
Set-TextValue "D15" "59.903.43"
$ws.Range("E15").Value = "  +3.37%  "

# Row 16
Set-TextValue "D16" "0.0000138"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17
Set-TextValue "D17" "2.419.05"
$ws.Range("E17").Value = "  +3.85%  "

# Row 18
Set-TextValue "D18" "11.30"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19
Set-TextValue "D19" "4.39"
$ws.Range("E19").Value = "  +1.89%  "

# Row 20
Set-TextValue "D20" "330.64"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
Set-TextValue "D21" "6.70"
$ws.Range("E21").Value = "  -3.01%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
Set-TextValue "D23" "66.09"
$ws.Range("E23").Value = "  +3.62%  "

# Row 24
$ws.Range("E24").Value = "  +1.44%  "

# Row 25
Set-TextValue "D25" "8.88"
$ws.Range("E25").Value = "  +7.51%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
Set-TextValue "D27" "1.38"
$ws.Range("E27").Value = "  +4.91%  "

# Row 28
Set-TextValue "D28" "0.0₃0780"
$ws.Range("E28").Value = "  +5.67%  "

# Row 29
$ws.Range("E29").Value = "  +0.53%  "

# Row 30
Set-TextValue "D30" "170.38"
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
Set-TextValue "D31" "6.13"
$ws.Range("E31").Value = "  -0.45%  "

# Row 32
Set-TextValue "D32" "18.69"
$ws.Range("E32").Value = "  +1.74%  "

# Row 33
Set-TextValue "D33" "1.02"
$ws.Range("E33").Value = "  +1.61%  "

# Row 34
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
Set-TextValue "D35" "1.30"
$ws.Range("E35").Value = "  +4.71%  "

# Row 36
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
Set-TextValue "D37" "4.22"
$ws.Range("E37").Value = "  +1.48%  "

# Row 38
Set-TextValue "D38" "1.61"
$ws.Range("E38").Value = "  +1.07%  "

# Row 39
Set-TextValue "D39" "39.64"
$ws.Range("E39").Value = "  -1.64%  "

# Row 40
Set-TextValue "D40" "0.411"
$ws.Range("E40").Value = "  -4.52%  "

# Row 41
Set-TextValue "D41" "314.11"
$ws.Range("E41").Value = "  +9.36%  "

# Row 42
Set-TextValue "D42" "3.69"
$ws.Range("E42").Value = "  +0.72%  "

# Row 43
Set-TextValue "D43" "139.00"
$ws.Range("E43").Value = "  -1.57%  "

# Row 44
Set-TextValue "D44" "0.0972"
$ws.Range("E44").Value = "  +1.47%  "

# Row 45
Set-TextValue "D45" "0.0520"
$ws.Range("E45").Value = "  +0.96%  "

# Row 46
Set-TextValue "D46" "19.54"
$ws.Range("E46").Value = "  +5.34%  "

# Row 47
Set-TextValue "D47" "0.580"
$ws.Range("E47").Value = "  +2.48%  "

# Row 48
Set-TextValue "D48" "0.0225"
$ws.Range("E48").Value = "  +1.46%  "

# Row 49
Set-TextValue "D49" "0.390"
$ws.Range("E49").Value = "  -3.01%  "

# Row 50
Set-TextValue "D50" "17.64"
$ws.Range("E50").Value = "  +1.33%  "

# Row 51
Set-TextValue "D51" "11.05"
$ws.Range("E51").Value = "  +0.34%  "
